$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.643.04'
$ws.Range('E2').Value = '  +1.11%  '
$ws.Range('D3').Value = '3.045.84'
$ws.Range('E3').Value = '  +2.97%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '384.46'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.09%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '102.91'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +0.71%  '
$ws.Range('E7').Value = '  -0.17%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.586'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -0.85%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '36.95'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +1.19%  '
$ws.Range('E11').Value = '  +0.20%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0864'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +1.15%  '
$ws.Range('D13').Value = '3.522.84'
$ws.Range('E13').Value = '  +2.82%  '
$ws.Range('E14').Value = '  +2.24%  '
$ws.Range('E15').Value = '  -0.47%  '
$ws.Range('D16').Value = '3.059.45'
$ws.Range('E16').Value = '  +3.28%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.979'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -1.77%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '10.59'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -8.03%  '
$ws.Range('D19').Value = '51.680.36'
$ws.Range('E19').Value = '  +1.02%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '3.11'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -0.55%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '12.48'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +0.74%  '
$ws.Range('D22').Value = '0.0₃0965'
$ws.Range('E22').Value = '  +0.55%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '70.03'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -0.13%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '267.26'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +0.23%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '3.17'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -4.09%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '8.42'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +6.97%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '7.38'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +2.70%  '
$ws.Range('E28').Value = '  +4.16%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '26.43'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +2.29%  '
$ws.Range('E30').Value = '  +0.02%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.108'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -2.64%  '
$ws.Range('E32').Value = '  -0.28%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '34.14'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -0.59%  '
$ws.Range('E34').Value = '  -2.43%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '50.56'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -1.02%  '
$ws.Range('E36').Value = '  +2.30%  '
$ws.Range('E37').Value = '  -0.14%  '
$ws.Range('E38').Value = '  +4.50%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.287'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +6.08%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '17.05'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +3.12%  '
$ws.Range('E41').Value = '  +2.27%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '128.73'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +2.71%  '
$ws.Range('E43').Value = '  -0.23%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.54'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +1.43%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '3.68'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +4.09%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '21.77'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +1.26%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.50'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +5.46%  '
$ws.Range('E48').Value = '  +3.24%  '
$ws.Range('D49').Value = '2.035.73'
$ws.Range('E49').Value = '  -0.77%  '
$ws.Range('D50').Value = '3.347.57'
$ws.Range('E50').Value = '  +2.91%  '
$ws.Range('E51').Value = '  +7.40%  '
